$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new "amelioration" rows describing the div/row overflow fix
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "Modification de la div row padding a 0px pour ne pas que la div dépasse"
$ws.Cells.Item(10, 3).Value = "fait"

$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "suppresion du css dans l'html"
$ws.Cells.Item(11, 3).Value = "fait"

$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "Modification des différentes parti du htlm"
$ws.Cells.Item(12, 3).Value = "fait"

# Column B needs to widen to fit the new, longer descriptions
$ws.Columns.Item(2).ColumnWidth = 65

# Move/restore the active selection to A13, matching the author's final cursor spot
$ws.Range("A13").Select()
